# Fixed a stupid bug with droplets walking - recalibration of droplet 18
# motor step values (the first "Droplet 18" table, rows 4-11).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Step 0 (row 4): only the final direction value changes.
$ws.Range("E4").Value = 207

# Steps 1-5 (rows 5-9): the per-step calibration offsets are no longer
# valid, so wipe them back to blank - only the step index in column B stays.
$ws.Range("C5:E9").ClearContents()

# Step 6 (row 10): tweak D/E, C stays the same (-1).
$ws.Range("D10").Value = -5
$ws.Range("E10").Value = 40

# Step 7 (row 11): new C/D calibration values, E stays the same (-1).
$ws.Range("C11").Value = 69
$ws.Range("D11").Value = 57

# Move the active selection to reflect where the user ended up working.
$ws.Range("F11").Select() | Out-Null
